# Delete the record with Kayıt No 11149144 (row 706 on "Kayitlar",
# row 164 on "Merkez İlçe") — shifts all following rows up by one.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(706).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(164).Delete()
